# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500adfae01c9a5dd7ef65e90accc96781b5c
# Rebrand from "ibm.com/Alvearie" to "linuxforhealth.org/LinuxForHealth", bump version/date,
# and correct the ele-1/ext-1 constraint placement on the Elements sheet.

$wb = $excel.ActiveWorkbook

# ----- Metadata sheet -----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/conversation-type"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ----- Elements sheet -----
$elements = $wb.Worksheets.Item("Elements")

# The generic ele-1 constraint no longer applies to the root "Extension" row (row 2);
# it now only shows up lower in the table (Extension.extension / Extension.value[x] rows).
$elements.Range("AI2").Value = ""

# Re-point the extension's Fixed Value / Binding Value Set at the new linuxforhealth.org host.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/conversation-type"
$elements.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/eng-conversation-type"

# Column Y widened (bestFit) to accommodate the longer linuxforhealth.org URL.
$elements.Columns.Item(25).ColumnWidth = 60.65
